$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 51 (formatting copied from row 51 automatically)
$ws.Rows("52:53").Insert()
$ws.Rows("52:53").RowHeight = 17

# New "Cycle Peak Labelling" / "Cycle Peak Symbols" shortcut rows
$ws.Range("A52").Value = "Cycle Peak Labelling"
$ws.Range("B52").Value = "PL"
$ws.Range("C52").Value = "PS"

$ws.Range("A53").Value = "Cycle Peak Symbols"
$ws.Range("B53").Value = "PS"
$ws.Range("C53").Value = "PL"

# Update the visible scroll position / selection
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C54").Select()

Write-Output "done"
